# "created wrapper and demo"
# Adds a new "cash" resource row (row 24) to the resources table, mirroring
# the formatting of the row above it, and leaves the selection on C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry over the existing row's formatting (font/style) onto the new row
# before writing the new values, so the new cells match the table's look.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)

$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)

$ws.Range("A24").Value = "cash"
$ws.Range("B24").Value = 0

$ws.Range("C23").Select()
